$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Menu Category" -> "Items"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Items"

$ws1.Range("A1").Value = "Category"
$ws1.Range("B1").Value = "Menu Item"
$ws1.Range("C1").Value = "Description"
$ws1.Range("D1").Value = "Costs"
$ws1.Range("E1").Value = "Option Groups"

# Copy the existing bold header style onto the newly used header cells
$ws1.Range("B1").Copy()
$ws1.Range("C1:E1").PasteSpecial(-4122)

$ws1.Range("A2").Value = "cold drinks"
$ws1.Range("B2").Value = "Menu Item"
$ws1.Range("C2").Value = "none"
$ws1.Range("D2").Value = 3.5
$ws1.Range("E2").Value = "cold options, can drinks range, hot options"

$ws1.Range("A3").Value = "hot drinks"
$ws1.Range("B3").Value = "Menu Item"
$ws1.Range("C3").Value = "none"
$ws1.Range("D3").Value = 3.5
$ws1.Range("E3").Value = "hot options"

# ---------------------------------------------------------------------------
# Sheet 2: "Linked Category" -> "Option Group"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Option Group"

$ws2.Range("A1").Value = "Option Group"
$ws2.Range("B1").Value = "Single"
$ws2.Range("C1").Value = "Mandatory"

$ws2.Range("B1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)

$ws2.Range("A2").Value = "cold options"
$ws2.Range("B2").Value = $True
$ws2.Range("C2").Value = $True

$ws2.Range("A3").Value = "can drinks range"
$ws2.Range("B3").Value = $True
$ws2.Range("C3").Value = $True

$ws2.Range("A4").Value = "hot options"
$ws2.Range("B4").Value = $True
$ws2.Range("C4").Value = $True

# ---------------------------------------------------------------------------
# Sheet 3: "Option Groups" -> "Options"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Options"

$ws3.Range("A2:B4").ClearContents()

$ws3.Range("A1").Value = "Option Group"
$ws3.Range("B1").Value = "Option"
$ws3.Range("C1").Value = "Cost"

$ws3.Range("B1").Copy()
$ws3.Range("C1").PasteSpecial(-4122)
